$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.755.42"
$ws.Range("E2").Value = "  -1.24%  "
$ws.Range("D3").Value = "3.370.90"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.86"
$ws.Range("E5").Value = "  -0.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.91"
$ws.Range("E6").Value = "  -0.64%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "3.369.32"
$ws.Range("E8").Value = "  -0.53%  "
$ws.Range("E9").Value = "  -1.18%  "
$ws.Range("E10").Value = "  +1.38%  "
$ws.Range("E11").Value = "  -3.29%  "
$ws.Range("E12").Value = "  -2.96%  "
$ws.Range("D13").Value = "3.946.15"
$ws.Range("E13").Value = "  -0.45%  "
$ws.Range("E14").Value = "  -0.71%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.96"
$ws.Range("E15").Value = "  +0.44%  "
$ws.Range("D16").Value = "3.371.03"
$ws.Range("E16").Value = "  -0.41%  "
$ws.Range("E17").Value = "  -3.89%  "
$ws.Range("D18").Value = "60.840.19"
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.76"
$ws.Range("E19").Value = "  -2.76%  "
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.80"
$ws.Range("E20").Value = "  -1.26%  "
$ws.Range("E21").Value = "  -2.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "371.53"
$ws.Range("E22").Value = "  -1.51%  "
$ws.Range("D23").Value = "3.507.11"
$ws.Range("E23").Value = "  -0.56%  "
$ws.Range("E24").Value = "  -1.93%  "
$ws.Range("E25").Value = "  +0.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "70.62"
$ws.Range("E26").Value = "  -0.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000123"
$ws.Range("E27").Value = "  -2.85%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.175"
$ws.Range("E28").Value = "  +7.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.59"
$ws.Range("E29").Value = "  -8.01%  "
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.32"
$ws.Range("E31").Value = "  -2.59%  "
$ws.Range("E32").Value = "  -2.87%  "
$ws.Range("E33").Value = "  -2.38%  "
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.27"
$ws.Range("E35").Value = "  -0.80%  "
$ws.Range("E36").Value = "  -3.83%  "
$ws.Range("E37").Value = "  -1.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.76"
$ws.Range("E38").Value = "  -1.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "164.49"
$ws.Range("E39").Value = "  -0.68%  "
$ws.Range("E40").Value = "  -2.72%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "25.48"
$ws.Range("E41").Value = "  +2.73%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.73"
$ws.Range("E43").Value = "  +0.63%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.769"
$ws.Range("E44").Value = "  -1.18%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.84"
$ws.Range("E45").Value = "  +0.90%  "
$ws.Range("E46").Value = "  -2.22%  "
$ws.Range("E47").Value = "  -6.30%  "
$ws.Range("D48").Value = "2.508.74"
$ws.Range("E48").Value = "  +6.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.50"
$ws.Range("E49").Value = "  +3.17%  "
$ws.Range("E50").Value = "  -1.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.40"
$ws.Range("E51").Value = "  +1.14%  "
